$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# Add a new "account number" column G, derived from the phone number in column F
# (replace the leading "08" digits with "00").
for ($r = 1; $r -le 10; $r++) {
    $ws.Cells.Item($r, 6).Copy()
    $ws.Cells.Item($r, 7).PasteSpecial($xlPasteFormats)

    $phone = $ws.Cells.Item($r, 6).Value()
    $suffix = $phone.Substring(3)
    $acct = "00$suffix"
    $ws.Cells.Item($r, 7).Value = $acct
}

$excel.CutCopyMode = $false
